$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the series-id header in B1 from "dc_real" to the FRED/Z.1 series
# code "BOGZ1FU315000005A".
$ws.Cells.Item(1, 2).Value2 = "BOGZ1FU315000005A"

# The values in column B (rows 2-50) were stored in raw dollars; rescale
# them to millions of dollars (divide by 1,000,000) to match the fixed data.
for ($r = 2; $r -le 50; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 / 1000000
}
